# Generate Report for Handback
# Update the generated/handoff/handback timestamps to reflect the latest run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row
$wsOverview.Range("G2").Value = "2016-08-20 09:10:12"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
$wsZhCn.Range("H2").Value = "2016-08-20 09:10:07"
$wsZhCn.Range("K2").Value = "2016-08-20 09:10:25"

# de-de sheet: "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-08-20 09:10:32"
